$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update collaborator counts for PR#0
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 10

# Update collaborator counts for PR#1
$ws.Range("D3").Value = 5

# Update collaborator counts for PR#2
$ws.Range("E4").Value = 4

# Update collaborator counts for PR#3
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 6

# Update collaborator counts for PR#6
$ws.Range("D8").Value = 5

# Update the active cell selection
$ws.Range("E4").Select()
